# The workbook was re-downloaded/saved from the web app; the only actual
# data change is on the "Аксесуари" sheet (the first sheet / sheetId=1),
# where the header cell A1 was re-written from lowercase "id" to "Id".
$wb = $excel.ActiveWorkbook

$wsAccessories = $wb.Worksheets.Item(1)
$wsAccessories.Range("A1").Value = "Id"
